# "fixed bugs with classifier"
# Appends three newly-classified cards to the bottom of the price_sheet
# table (rows 282-284), matching the formatting of the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("price_sheet")

# Carry the formatting (fonts/fills/borders/number format) of the last
# existing data row down into the three new rows before filling values.
$ws.Range("A281:D281").Copy() | Out-Null
$ws.Range("A282:D284").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 282: Worm Token, Common, $0.07, GRN
$ws.Range("A282").Value = "Worm Token"
$ws.Range("B282").Value = "C"
$ws.Range("C282").Value = 0.07
$ws.Range("D282").Value = "GRN"

# Row 283: Warrior Token, Common, $0.04, GRN
$ws.Range("A283").Value = "Warrior Token"
$ws.Range("B283").Value = "C"
$ws.Range("C283").Value = 0.04
$ws.Range("D283").Value = "GRN"

# Row 284: Sphinx Insight, Uncommon, $0.07, GRN
$ws.Range("A284").Value = "Sphinx Insight"
$ws.Range("B284").Value = "U"
$ws.Range("C284").Value = 0.07
$ws.Range("D284").Value = "GRN"

# Leave the view scrolled/selected near the newly-added rows, matching
# where the author was working when the file was saved.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 263
$win.ScrollColumn = 1
$ws.Range("G280").Select() | Out-Null
